$d = $word.ActiveDocument

# --- Paragraph 3: "Compare and contrast bioinformatics and computational biology" ---
# -> bold "bioinformatics" and "computational biology"
$p1 = $d.Paragraphs.Item(3)
$r1 = $p1.Range
$s1 = $r1.Start
$t1 = $r1.Text

$i1 = $t1.IndexOf("bioinformatics")
$sub1 = $d.Range($s1 + $i1, $s1 + $i1 + "bioinformatics".Length)
$sub1.Font.Bold = 1

$i2 = $t1.IndexOf("computational biology")
$sub2 = $d.Range($s1 + $i2, $s1 + $i2 + "computational biology".Length)
$sub2.Font.Bold = 1

# --- Paragraph 4: "Compare and contrast computational biology/bioinformatics with other biological disciplines" ---
# -> bold "computational biology/bioinformatics" and remove the trailing _GoBack bookmark
$p2 = $d.Paragraphs.Item(4)
$r2 = $p2.Range
$s2 = $r2.Start
$t2 = $r2.Text

$i3 = $t2.IndexOf("computational biology/bioinformatics")
$sub3 = $d.Range($s2 + $i3, $s2 + $i3 + "computational biology/bioinformatics".Length)
$sub3.Font.Bold = 1

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# --- Paragraph 6: "Synthesize your knowledge about biology to identify types of "big data" in biology" ---
# -> bold "big data" and wrap it with a new _GoBack bookmark
$p3 = $d.Paragraphs.Item(6)
$r3 = $p3.Range
$s3 = $r3.Start
$t3 = $r3.Text

$i4 = $t3.IndexOf("big data")
$sub4 = $d.Range($s3 + $i4, $s3 + $i4 + "big data".Length)
$sub4.Font.Bold = 1
$d.Bookmarks.Add("_GoBack", $sub4)
